$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the literal 0 values in G10:G113 while keeping their number formatting/styles.
$ws.Range("G10:G113").ClearContents()

# 2. Add a thin gray top border to row 114 (B114:G114) - visually separates the
#    detail rows above from the totals section below. (Color is theme 1 / "Text 1"
#    tinted +35%, i.e. RGB #595959 -> decimal 5855577.)
$sepRange = $ws.Range("B114:G114")
$topBorder = $sepRange.Borders.Item(8)
$topBorder.Color = 5855577
$topBorder.Weight = 2
$topBorder.LineStyle = 1

# 3. Update the sheet view: select G10:G113 (this also clears the previously
#    scrolled-down top-left cell, resetting the view back to the top).
$ws.Activate()
$ws.Range("G10:G113").Select()
